$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheetPr/outline + page-margin settings used by the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 3. Copy cell formatting from the existing sheets so the new sheet ---
#        matches the look (bold/centered header row, date-formatted A column)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A82").PasteSpecial(-4122)

# --- 4. Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 5. Data rows (2 - 82) ---
$data = New-Object 'object[,]' 81,4
$data[0,0] = 44934.99999999999
$data[0,1] = 103
$data[0,2] = -95.47299678968157
$data[0,3] = 309.7796319917326
$data[1,0] = 44941.99999999999
$data[1,1] = 104
$data[1,2] = -105.9340868629948
$data[1,3] = 304.5566427850823
$data[2,0] = 44948.99999999999
$data[2,1] = 105
$data[2,2] = -99.29509449727175
$data[2,3] = 298.0605601928618
$data[3,0] = 44955.99999999999
$data[3,1] = 106
$data[3,2] = -94.30919062920204
$data[3,3] = 301.8528278711614
$data[4,0] = 44962.99999999999
$data[4,1] = 107
$data[4,2] = -102.8527822472914
$data[4,3] = 316.2915951844864
$data[5,0] = 44969.99999999999
$data[5,1] = 108
$data[5,2] = -104.288529926833
$data[5,3] = 309.0009816623187
$data[6,0] = 44976.99999999999
$data[6,1] = 109
$data[6,2] = -108.667129703273
$data[6,3] = 306.0936592710922
$data[7,0] = 44983.99999999999
$data[7,1] = 110
$data[7,2] = -99.95933657439707
$data[7,3] = 302.3043418822893
$data[8,0] = 44990.99999999999
$data[8,1] = 111
$data[8,2] = -88.72576413492924
$data[8,3] = 314.6045896457301
$data[9,0] = 44997.99999999999
$data[9,1] = 112
$data[9,2] = -81.90237287367142
$data[9,3] = 343.9281914118201
$data[10,0] = 45004.99999999999
$data[10,1] = 113
$data[10,2] = -86.86131405519254
$data[10,3] = 312.1854474042227
$data[11,0] = 45011.99999999999
$data[11,1] = 114
$data[11,2] = -74.41254932613788
$data[11,3] = 332.6973055828349
$data[12,0] = 45018.99999999999
$data[12,1] = 115
$data[12,2] = -94.24895686435001
$data[12,3] = 315.949312779964
$data[13,0] = 45025.99999999999
$data[13,1] = 116
$data[13,2] = -90.11836482171468
$data[13,3] = 318.5117015534476
$data[14,0] = 45032.99999999999
$data[14,1] = 117
$data[14,2] = -86.65166580006385
$data[14,3] = 327.0132404319524
$data[15,0] = 45039.99999999999
$data[15,1] = 118
$data[15,2] = -91.17250495702025
$data[15,3] = 321.2416053927175
$data[16,0] = 45046.99999999999
$data[16,1] = 119
$data[16,2] = -87.03100602947009
$data[16,3] = 321.5509265748955
$data[17,0] = 45053.99999999999
$data[17,1] = 120
$data[17,2] = -76.14841601779163
$data[17,3] = 323.276933095053
$data[18,0] = 45060.99999999999
$data[18,1] = 121
$data[18,2] = -69.07865902442285
$data[18,3] = 330.5834663425001
$data[19,0] = 45067.99999999999
$data[19,1] = 122
$data[19,2] = -91.9621949379585
$data[19,3] = 317.9161082825771
$data[20,0] = 45074.99999999999
$data[20,1] = 123
$data[20,2] = -80.95406825763486
$data[20,3] = 341.6913028132257
$data[21,0] = 45081.99999999999
$data[21,1] = 124
$data[21,2] = -75.27440768150628
$data[21,3] = 313.2496168946919
$data[22,0] = 45088.99999999999
$data[22,1] = 125
$data[22,2] = -74.01312453192001
$data[22,3] = 351.7355778913904
$data[23,0] = 45095.99999999999
$data[23,1] = 126
$data[23,2] = -89.61674258015
$data[23,3] = 337.8617524680021
$data[24,0] = 45102.99999999999
$data[24,1] = 127
$data[24,2] = -73.25558099769444
$data[24,3] = 332.9975019058663
$data[25,0] = 45109.99999999999
$data[25,1] = 128
$data[25,2] = -73.46005132914149
$data[25,3] = 329.5841582168534
$data[26,0] = 45116.99999999999
$data[26,1] = 129
$data[26,2] = -82.45088305137365
$data[26,3] = 329.6668903576673
$data[27,0] = 45123.99999999999
$data[27,1] = 130
$data[27,2] = -90.69553743978719
$data[27,3] = 334.2450336236761
$data[28,0] = 45130.99999999999
$data[28,1] = 131
$data[28,2] = -73.3474692016102
$data[28,3] = 322.7908782732153
$data[29,0] = 45137.99999999999
$data[29,1] = 132
$data[29,2] = -61.98082539341808
$data[29,3] = 331.546247335831
$data[30,0] = 45144.99999999999
$data[30,1] = 133
$data[30,2] = -85.95505517177291
$data[30,3] = 337.2980016116566
$data[31,0] = 45151.99999999999
$data[31,1] = 134
$data[31,2] = -72.29811229734065
$data[31,3] = 351.4646538318501
$data[32,0] = 45158.99999999999
$data[32,1] = 135
$data[32,2] = -69.67840167889179
$data[32,3] = 344.6469700755173
$data[33,0] = 45165.99999999999
$data[33,1] = 135
$data[33,2] = -72.36923994044118
$data[33,3] = 332.5244963598668
$data[34,0] = 45172.99999999999
$data[34,1] = 136
$data[34,2] = -57.86361279106836
$data[34,3] = 343.3165771141298
$data[35,0] = 45179.99999999999
$data[35,1] = 137
$data[35,2] = -49.04903143556663
$data[35,3] = 343.9000805698813
$data[36,0] = 45186.99999999999
$data[36,1] = 138
$data[36,2] = -76.73872142865524
$data[36,3] = 349.3168211603813
$data[37,0] = 45193.99999999999
$data[37,1] = 139
$data[37,2] = -65.48075855508105
$data[37,3] = 319.7869873284597
$data[38,0] = 45200.99999999999
$data[38,1] = 140
$data[38,2] = -67.99601879275194
$data[38,3] = 335.8438037280418
$data[39,0] = 45207.99999999999
$data[39,1] = 141
$data[39,2] = -70.09646841462038
$data[39,3] = 347.9236716611505
$data[40,0] = 45214.99999999999
$data[40,1] = 142
$data[40,2] = -61.47183950206269
$data[40,3] = 347.5501376426461
$data[41,0] = 45221.99999999999
$data[41,1] = 143
$data[41,2] = -66.31558910134727
$data[41,3] = 352.1849718809729
$data[42,0] = 45228.99999999999
$data[42,1] = 144
$data[42,2] = -62.05878628176714
$data[42,3] = 354.8832491138523
$data[43,0] = 45235.99999999999
$data[43,1] = 145
$data[43,2] = -54.80146935182574
$data[43,3] = 341.9798697889512
$data[44,0] = 45242.99999999999
$data[44,1] = 146
$data[44,2] = -67.9299379963913
$data[44,3] = 342.0233249384747
$data[45,0] = 45249.99999999999
$data[45,1] = 147
$data[45,2] = -64.95840108279239
$data[45,3] = 349.1817132473888
$data[46,0] = 45256.99999999999
$data[46,1] = 148
$data[46,2] = -57.83176224590582
$data[46,3] = 351.2793662762562
$data[47,0] = 45263.99999999999
$data[47,1] = 149
$data[47,2] = -42.4448734060574
$data[47,3] = 361.0191565329225
$data[48,0] = 45270.99999999999
$data[48,1] = 150
$data[48,2] = -51.69527084233162
$data[48,3] = 366.5023703706864
$data[49,0] = 45277.99999999999
$data[49,1] = 151
$data[49,2] = -43.8771999948304
$data[49,3] = 357.4712857900308
$data[50,0] = 45298.99999999999
$data[50,1] = 154
$data[50,2] = -56.26878415529787
$data[50,3] = 361.3929193838338
$data[51,0] = 45305.99999999999
$data[51,1] = 155
$data[51,2] = -49.27273541915288
$data[51,3] = 355.4080611505872
$data[52,0] = 45312.99999999999
$data[52,1] = 156
$data[52,2] = -52.97076607028016
$data[52,3] = 359.3476432724116
$data[53,0] = 45319.99999999999
$data[53,1] = 157
$data[53,2] = -56.65488371648296
$data[53,3] = 345.1358278979941
$data[54,0] = 45326.99999999999
$data[54,1] = 158
$data[54,2] = -52.21465311595235
$data[54,3] = 356.9877718061797
$data[55,0] = 45333.99999999999
$data[55,1] = 159
$data[55,2] = -49.50461914972409
$data[55,3] = 365.61261203396
$data[56,0] = 45347.99999999999
$data[56,1] = 161
$data[56,2] = -43.32451646804225
$data[56,3] = 352.4190024979924
$data[57,0] = 45354.99999999999
$data[57,1] = 162
$data[57,2] = -44.31257973520594
$data[57,3] = 373.2076899066022
$data[58,0] = 45361.99999999999
$data[58,1] = 163
$data[58,2] = -61.26586197776274
$data[58,3] = 369.6859949166515
$data[59,0] = 45382.99999999999
$data[59,1] = 166
$data[59,2] = -37.84707604453478
$data[59,3] = 383.6112348371003
$data[60,0] = 45389.99999999999
$data[60,1] = 167
$data[60,2] = -22.33035826944227
$data[60,3] = 368.5503791648592
$data[61,0] = 45459.99999999999
$data[61,1] = 177
$data[61,2] = -30.64546399363025
$data[61,3] = 371.6937439456383
$data[62,0] = 45466.99999999999
$data[62,1] = 178
$data[62,2] = -33.75795895894797
$data[62,3] = 386.5957151286267
$data[63,0] = 45473.99999999999
$data[63,1] = 179
$data[63,2] = -33.10822463700374
$data[63,3] = 394.4044280432256
$data[64,0] = 45480.99999999999
$data[64,1] = 180
$data[64,2] = -22.30166073738711
$data[64,3] = 374.8819242801237
$data[65,0] = 45487.99999999999
$data[65,1] = 180
$data[65,2] = -29.00170988221141
$data[65,3] = 393.3916781114776
$data[66,0] = 45494.99999999999
$data[66,1] = 181
$data[66,2] = -28.5902675315464
$data[66,3] = 384.207506201367
$data[67,0] = 45515.99999999999
$data[67,1] = 184
$data[67,2] = -16.53623353456505
$data[67,3] = 376.0139122614681
$data[68,0] = 45529.99999999999
$data[68,1] = 186
$data[68,2] = -18.22374318645742
$data[68,3] = 383.5050146746414
$data[69,0] = 45536.99999999999
$data[69,1] = 187
$data[69,2] = 10.02023330853343
$data[69,3] = 390.5999377905579
$data[70,0] = 45543.99999999999
$data[70,1] = 188
$data[70,2] = -10.59790395567223
$data[70,3] = 387.8092264818146
$data[71,0] = 45550.99999999999
$data[71,1] = 189
$data[71,2] = -18.6337406127917
$data[71,3] = 375.9964168851479
$data[72,0] = 45641.99999999999
$data[72,1] = 202
$data[72,2] = 9.151966084151855
$data[72,3] = 401.650875455896
$data[73,0] = 45648.99999999999
$data[73,1] = 203
$data[73,2] = -5.781378773519132
$data[73,3] = 407.709348289282
$data[74,0] = 45655.99999999999
$data[74,1] = 204
$data[74,2] = -7.376282602451892
$data[74,3] = 417.6074075649153
$data[75,0] = 45662.99999999999
$data[75,1] = 205
$data[75,2] = 12.64273688442991
$data[75,3] = 397.5231154129514
$data[76,0] = 45669.99999999999
$data[76,1] = 206
$data[76,2] = -4.989390785671088
$data[76,3] = 412.0624766667578
$data[77,0] = 45676.99999999999
$data[77,1] = 207
$data[77,2] = 8.644752465761377
$data[77,3] = 405.7621240480015
$data[78,0] = 45683.99999999999
$data[78,1] = 208
$data[78,2] = -4.013108285107387
$data[78,3] = 419.3202451813526
$data[79,0] = 45690.99999999999
$data[79,1] = 209
$data[79,2] = 15.05619014193309
$data[79,3] = 399.445539642326
$data[80,0] = 45697.99999999999
$data[80,1] = 210
$data[80,2] = 8.783799914522987
$data[80,3] = 418.5601552716823

$wsForecast.Range("A2:D82").Value = $data
